$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = "'22.388.97"
$ws.Range('D2').Style = 'Normal'
$ws.Range('E2').Value = '  -0.45%  '
$ws.Range('D3').Value = "'1.573.97"
$ws.Range('D3').Style = 'Normal'
$ws.Range('E3').Value = '  +0.00%  '
$ws.Range('E4').Value = '  -0.35%  '
$ws.Range('E5').Value = '  -0.36%  '
$ws.Range('D6').Value = "'290.92"
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  -0.71%  '
$ws.Range('D7').Value = "'0.3760"
$ws.Range('D7').Style = 'Normal'
$ws.Range('E7').Value = '  +2.18%  '
$ws.Range('D8').Value = "'50.10"
$ws.Range('D8').Style = 'Normal'
$ws.Range('E8').Value = '  +1.63%  '
$ws.Range('D9').Value = "'0.3416"
$ws.Range('D9').Style = 'Normal'
$ws.Range('E9').Value = '  +2.70%  '
$ws.Range('D10').Value = "'1.166"
$ws.Range('D10').Style = 'Normal'
$ws.Range('E10').Value = '  +0.32%  '
$ws.Range('D11').Value = "'0.07681"
$ws.Range('D11').Style = 'Normal'
$ws.Range('E11').Value = '  +1.68%  '
$ws.Range('E12').Value = '  -0.30%  '
$ws.Range('D13').Value = "'21.37"
$ws.Range('D13').Style = 'Normal'
$ws.Range('E13').Value = '  +1.58%  '
$ws.Range('D14').Value = "'5.988"
$ws.Range('D14').Style = 'Normal'
$ws.Range('E14').Value = '  -0.97%  '
$ws.Range('D15').Value = "'6.931"
$ws.Range('D15').Style = 'Normal'
$ws.Range('E15').Value = '  +0.98%  '
$ws.Range('E16').Value = '  +0.81%  '
$ws.Range('D17').Value = "'1.573.27"
$ws.Range('D17').Style = 'Normal'
$ws.Range('E17').Value = '  +0.24%  '
$ws.Range('D18').Value = "'90.36"
$ws.Range('D18').Style = 'Normal'
$ws.Range('E18').Value = '  +0.90%  '
$ws.Range('D19').Value = "'0.06722"
$ws.Range('D19').Style = 'Normal'
$ws.Range('E19').Value = '  -0.46%  '
$ws.Range('E20').Value = '  -0.32%  '
$ws.Range('D21').Value = "'16.77"
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = '  +2.85%  '
$ws.Range('E22').Value = '  +0.11%  '
$ws.Range('D23').Value = "'0.5273"
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = '  -4.27%  '
$ws.Range('E24').Value = '  +1.59%  '
$ws.Range('D25').Value = "'22.392.73"
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = '  -0.49%  '
$ws.Range('E26').Value = '  +1.40%  '
$ws.Range('D27').Value = "'2.781"
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = '  -4.29%  '
$ws.Range('D28').Value = "'20.31"
$ws.Range('D28').Style = 'Normal'
$ws.Range('E28').Value = '  +3.23%  '
$ws.Range('D29').Value = "'144.70"
$ws.Range('D29').Style = 'Normal'
$ws.Range('E29').Value = '  -0.81%  '
$ws.Range('D30').Value = "'5.072"
$ws.Range('D30').Style = 'Normal'
$ws.Range('E30').Value = '  +2.74%  '
$ws.Range('E31').Value = '  +0.93%  '
$ws.Range('D32').Value = "'1.746.33"
$ws.Range('D32').Style = 'Normal'
$ws.Range('E32').Value = '  -0.17%  '
$ws.Range('D33').Value = "'1.024"
$ws.Range('D33').Style = 'Normal'
$ws.Range('E33').Value = '  +8.87%  '
$ws.Range('D34').Value = "'6.249"
$ws.Range('D34').Style = 'Normal'
$ws.Range('E34').Value = '  +0.28%  '
$ws.Range('D35').Value = "'2.024"
$ws.Range('D35').Style = 'Normal'
$ws.Range('E35').Value = '  -0.21%  '
$ws.Range('E36').Value = '  -2.76%  '
$ws.Range('D37').Value = "'0.08522"
$ws.Range('D37').Style = 'Normal'
$ws.Range('E37').Value = '  -0.37%  '
$ws.Range('D38').Value = "'0.02547"
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').Value = '  +1.92%  '
$ws.Range('D39').Value = "'0.2326"
$ws.Range('D39').Style = 'Normal'
$ws.Range('E39').Value = '  +2.06%  '
$ws.Range('D40').Value = "'0.06539"
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value = '  +0.50%  '
$ws.Range('D41').Value = "'5.512"
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Value = '  +2.19%  '
$ws.Range('E42').Value = '  +2.63%  '
$ws.Range('D43').Value = "'11.68"
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = '  -0.11%  '
$ws.Range('D44').Value = "'0.6438"
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = '  +1.72%  '
$ws.Range('D45').Value = "'14.08"
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value = '  -2.64%  '
$ws.Range('E46').Value = '  -0.31%  '
$ws.Range('D47').Value = "'0.6023"
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = '  +1.81%  '
$ws.Range('D48').Value = "'3.783"
$ws.Range('D48').Style = 'Normal'
$ws.Range('E48').Value = '  +0.08%  '
$ws.Range('D49').Value = "'1.304"
$ws.Range('D49').Style = 'Normal'
$ws.Range('E49').Value = '  +11.59%  '
$ws.Range('D50').Value = "'2.098"
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = '  -0.06%  '
$ws.Range('D51').Value = "'125.03"
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = '  +2.21%  '
